$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.615.50"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.392.25"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.384"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "3.970.94"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "3.394.30"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "61.657.99"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.183"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "168.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").Value = "3.424.39"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0756"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("D45").Value = "2.477.17"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -6.33%  "
$ws.Range("E51").Value = "  -1.41%  "
